$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Copy the formatting of the last populated "normal" row (34) down onto the
# previously-empty rows 35-41 so the new cells pick up the existing shared
# style (border/fill/font) instead of creating a brand-new style entry.
$ws.Range("A34:C34").Copy()
$ws.Range("A35:C41").PasteSpecial(-4122)

# Row 42's A/C cells use that same recurring style, but B42 uses the
# special "Docs-Calibri" font style already present on a couple of other
# cells in the sheet (A28/A30/B30) - reuse that one as the base.
$ws.Range("A34").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("C34").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("A28").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 35
$ws.Range("A35").Value = "Automation_Folder_Inside_AutomationFolder"
$ws.Range("B35").Value = "Automation_Folder_Inside_AutomationFolder"
$ws.Range("C35").Value = 'New Folder created inside folder "Automation_Folder_ToCopyDashboardToIt"'

# Row 36
$ws.Range("A36").Value = "Automation_Dashboard_CopiedAndFoundInSearch"
$ws.Range("B36").Value = "Automation_Dashboard_CopiedAndFoundInSearch"
$ws.Range("C36").Value = "New Dashboard"

# Row 37
$ws.Range("A37").Value = "Automation_Dashboard_DuplicateError Copy"
$ws.Range("B37").Value = "Automation_Dashboard_DuplicateError Copy"
$ws.Range("C37").Value = 'New dashboard created inside folder "Automation_Folder_ToCopyDashboardToIt"'

# Row 38
$ws.Range("A38").Value = "Automation_Dashboard_DuplicateError"
$ws.Range("B38").Value = "Automation_Dashboard_DuplicateError"
$ws.Range("C38").Value = "New Dashboard"

# Row 39
$ws.Range("A39").Value = "Automation_Folder_Copy1"
$ws.Range("B39").Value = "Automation_Folder_Copy1"
$ws.Range("C39").Value = "New Folder"

# Row 40
$ws.Range("A40").Value = "Automation_Folder_Copy2"
$ws.Range("B40").Value = "Automation_Folder_Copy2"
$ws.Range("C40").Value = "New Folder"

# Row 41
$ws.Range("A41").Value = "Automation_Dashboard_Copy1"
$ws.Range("B41").Value = "Automation_Dashboard_Copy1"
$ws.Range("C41").Value = "New Dashboard"

# Row 42
$ws.Range("A42").Value = "Automation_Dashboard_Copy2"
$ws.Range("B42").Value = "Automation_Dashboard_Copy2"
$ws.Range("B42").Font.Name = "Docs-Calibri"
$ws.Range("C42").Value = 'New Dashboard exist in "Automation_Folder_Copy2"'
